$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (was D2=44875, M2=50, N2=16000, O2=16000, P2=16000, Q2="$/bandeja 10 kilos", S2=1600, T2=10)
# becomes the values previously in row 4.
$ws.Range("D2").Value = 44855
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/bandeja 5 kilos"
$ws.Range("S2").Value = 3000
$ws.Range("T2").Value = 5

# Row 4 (was D4=44855, M4=25, N4=15000, O4=15000, P4=15000, Q4="$/bandeja 5 kilos", S4=3000, T4=5)
# becomes the values previously in row 2.
$ws.Range("D4").Value = 44875
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1600
$ws.Range("T4").Value = 10
